# "create counter hotels page"
# Mark the hotel-related requirement rows as done by putting a 1 in column B
# (column B is the "counter" column; C already has 1 for every one of these
# requirement rows). B52 holds =SUM(B2:B51) and will recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

$rows = @(18, 19, 20, 22, 23, 24, 28, 29, 30, 31, 32, 49, 51)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 1
}

# Restore the on-screen selection to match where the author left off editing.
$ws.Range("F17").Select()
